$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) of the last existing data row (74) down into the
# two new rows (75, 76) so the new rows inherit the same cell styles
# (bold/border/center for column A, date number format for column E, etc.)
$ws.Range("A74:V74").Copy()
$ws.Range("A75:V76").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Row 75 ----
$r = 75
$ws.Cells.Item($r, 1).Value = 74
$ws.Cells.Item($r, 2).Value = "thailand"
$ws.Cells.Item($r, 3).Value = "thai-league-1"
$ws.Cells.Item($r, 4).Value = "2023-2024"
$ws.Cells.Item($r, 5).Value = 45236.52083333334
$ws.Cells.Item($r, 6).Value = "Trat FC"
$ws.Cells.Item($r, 7).Value = 0
$ws.Cells.Item($r, 8).Value = "Muang Thong Utd"
$ws.Cells.Item($r, 9).Value = 1
$ws.Cells.Item($r, 10).Value = 3.2
$ws.Cells.Item($r, 11).Value = "05/11/2023 16:38"
$ws.Cells.Item($r, 12).Value = 2.93
$ws.Cells.Item($r, 13).Value = "06/11/2023 12:21"
$ws.Cells.Item($r, 14).Value = 3.66
$ws.Cells.Item($r, 15).Value = "05/11/2023 16:38"
$ws.Cells.Item($r, 16).Value = 3.74
$ws.Cells.Item($r, 17).Value = "06/11/2023 12:21"
$ws.Cells.Item($r, 18).Value = 2.05
$ws.Cells.Item($r, 19).Value = "05/11/2023 16:38"
$ws.Cells.Item($r, 20).Value = 2.29
$ws.Cells.Item($r, 21).Value = "06/11/2023 12:21"
$ws.Cells.Item($r, 22).Value = "https://www.betexplorer.com/football/thailand/thai-league-1/trat-fc-muang-thong-utd/vZmmzlVO/"

# ---- Row 76 ----
$r = 76
$ws.Cells.Item($r, 1).Value = 75
$ws.Cells.Item($r, 2).Value = "thailand"
$ws.Cells.Item($r, 3).Value = "thai-league-1"
$ws.Cells.Item($r, 4).Value = "2023-2024"
$ws.Cells.Item($r, 5).Value = 45236.58333333334
$ws.Cells.Item($r, 6).Value = "Police Tero"
$ws.Cells.Item($r, 7).Value = 2
$ws.Cells.Item($r, 8).Value = "Lamphun Warrior"
$ws.Cells.Item($r, 9).Value = 2
$ws.Cells.Item($r, 10).Value = 1.87
$ws.Cells.Item($r, 11).Value = "30/10/2023 16:42"
$ws.Cells.Item($r, 12).Value = 2.16
$ws.Cells.Item($r, 13).Value = "06/11/2023 13:54"
$ws.Cells.Item($r, 14).Value = 3.71
$ws.Cells.Item($r, 15).Value = "30/10/2023 16:42"
$ws.Cells.Item($r, 16).Value = 3.72
$ws.Cells.Item($r, 17).Value = "06/11/2023 13:59"
$ws.Cells.Item($r, 18).Value = 3.7
$ws.Cells.Item($r, 19).Value = "30/10/2023 16:42"
$ws.Cells.Item($r, 20).Value = 3.19
$ws.Cells.Item($r, 21).Value = "06/11/2023 13:54"
$ws.Cells.Item($r, 22).Value = "https://www.betexplorer.com/football/thailand/thai-league-1/police-tero-lamphun-warrior/foOdfiGt/"
